# Add create and edit ticket api
#
# Two content changes to the "Data Dictionary" sheet:
#  1. The Ticket.Comments field description gains author_id/author.
#  2. The single "Log - [{User, Date, Action} ]" field on the Ticket History log table is
#     replaced by three separate fields (action, date, user), which pushes the Admin log
#     sub-table down two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the Ticket.Comments data-dictionary description.
$ws.Range("E15").Value = "Comments [ { author_id, author, name, date, content }]"

# 2) Relocate the "Admin log" block (K8:L12) down two rows, to K10:L14, to make room for
#    the new Ticket History log fields. Stage the move through an unused scratch range so
#    it behaves correctly even though source and destination overlap.
$ws.Range("K8:L12").Copy($ws.Range("Z8:AA12"))
$ws.Range("Z8:AA12").Copy($ws.Range("K10:L14"))
$ws.Range("Z8:AA12").Clear()

# Row 9 (K9:L9) held the old "ID" row of the Admin log block; once everything below it
# has shifted down, that leftover cell must be cleared out entirely.
$ws.Range("K9:L9").Clear()

# 3) Fill the vacated rows K6:L8 with the three new fields: action, date, user.
#    K6:L6 and K7:L7 get the regular "middle of box" formatting (same as K4:L4).
#    K8:L8 becomes the new last row of the box, so it gets the bottom-border formatting
#    (same family as E8:F8, the last row of the neighboring Projects box).
$ws.Range("K4:L4").Copy($ws.Range("K6:L6"))
$ws.Range("K4:L4").Copy($ws.Range("K7:L7"))
$ws.Range("E8:F8").Copy($ws.Range("K8:L8"))

$ws.Range("K8").Value = "user"
$ws.Range("K6").Value = "action"
$ws.Range("K7").Value = "date"
